$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "https://myweb3jobs.com/job/senior-product-marketing-manager/",
    "https://myweb3jobs.com/job/senior-solidity-engineer-6/",
    "https://myweb3jobs.com/job/associate-pr-manager/",
    "https://myweb3jobs.com/job/lead-front-end-engineer/",
    "https://myweb3jobs.com/job/content-growth-marketer/",
    "https://myweb3jobs.com/job/junior-compliance-analyst/",
    "https://myweb3jobs.com/job/senior-backend-engineer-20/",
    "https://myweb3jobs.com/job/full-stack-software-engineer-mid-level/",
    "https://myweb3jobs.com/job/user-researcher/",
    "https://myweb3jobs.com/job/senior-software-engineer-backend-4/",
    "https://myweb3jobs.com/job/community-manager-associate/",
    "https://myweb3jobs.com/job/product-marketer-2/",
    "https://myweb3jobs.com/job/senior-backend-engineer-trading-desk-engineering/",
    "https://myweb3jobs.com/job/video-editor-contractor/",
    "https://myweb3jobs.com/job/senior-network-engineer/"
)

$startRow = 7
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $values[$i]
}
